$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quote-prefix the date-format cells so they store literal text instead of numeric dates
$ws.Range("C2").Value = "'1/1/2017"
$ws.Range("D2").Value = "'11/15/2021"

# Re-apply the original date number format (entering a quote-prefixed value
# nudges Excel to create a brand new numFmt - restore the original m/d/yyyy)
$ws.Range("C2:D2").NumberFormat = "m/d/yyyy"

# Update the selection to match the recorded view state
$ws.Range("D6").Select()
